# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
#
# The worker-level "Periodo Mora" rows (B16:G40) are regrouped so that every
# period owed by a given worker is listed together (descending by period),
# two new workers (DAVINSON DARIO DIAZ LEDESMA / 1131106099 and MARCO
# AURELIO MORA GUZMAN / 73181246) pick up a full run of periods (2205..2111),
# and SEBASTIAN ENRIQUE VILLARREAL PERTUZ is extended back to period 2109.
# Row 16 (RODOLFO VALENZUELA GONZALEZ / 2105) is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => (Tipo Doc, N Doc, Nombre, Periodo, Valor Mora)  -- Salario Basico (G) stays 908526 throughout
$data = @(
    @(17, "CC", "1131106099", "DAVINSON DARIO DIAZ LEDESMA",          "2205", 26919),
    @(18, "CC", "1131106099", "DAVINSON DARIO DIAZ LEDESMA",          "2204", 35112),
    @(19, "CC", "1131106099", "DAVINSON DARIO DIAZ LEDESMA",          "2203", 36341),
    @(20, "CC", "1131106099", "DAVINSON DARIO DIAZ LEDESMA",          "2202", 36341),
    @(21, "CC", "1131106099", "DAVINSON DARIO DIAZ LEDESMA",          "2201", 36341),
    @(22, "CC", "1131106099", "DAVINSON DARIO DIAZ LEDESMA",          "2112", 36341),
    @(23, "CC", "1131106099", "DAVINSON DARIO DIAZ LEDESMA",          "2111", 36341),
    @(24, "CC", "73181246",   "MARCO AURELIO MORA GUZMAN",            "2205", 26919),
    @(25, "CC", "73181246",   "MARCO AURELIO MORA GUZMAN",            "2204", 35112),
    @(26, "CC", "73181246",   "MARCO AURELIO MORA GUZMAN",            "2203", 35112),
    @(27, "CC", "73181246",   "MARCO AURELIO MORA GUZMAN",            "2202", 35112),
    @(28, "CC", "73181246",   "MARCO AURELIO MORA GUZMAN",            "2201", 35112),
    @(29, "CC", "73181246",   "MARCO AURELIO MORA GUZMAN",            "2112", 35112),
    @(30, "CC", "73181246",   "MARCO AURELIO MORA GUZMAN",            "2111", 35112),
    @(31, "CC", "19752336",   "CRISTIAN ROCHA JIMENEZ",                "2105", 36341),
    @(32, "CC", "1143377371", "SEBASTIAN ENRIQUE VILLARREAL PERTUZ",  "2205", 27861),
    @(33, "CC", "1143377371", "SEBASTIAN ENRIQUE VILLARREAL PERTUZ",  "2204", 36341),
    @(34, "CC", "1143377371", "SEBASTIAN ENRIQUE VILLARREAL PERTUZ",  "2203", 36341),
    @(35, "CC", "1143377371", "SEBASTIAN ENRIQUE VILLARREAL PERTUZ",  "2202", 36341),
    @(36, "CC", "1143377371", "SEBASTIAN ENRIQUE VILLARREAL PERTUZ",  "2201", 36341),
    @(37, "CC", "1143377371", "SEBASTIAN ENRIQUE VILLARREAL PERTUZ",  "2112", 36341),
    @(38, "CC", "1143377371", "SEBASTIAN ENRIQUE VILLARREAL PERTUZ",  "2111", 36341),
    @(39, "CC", "1143377371", "SEBASTIAN ENRIQUE VILLARREAL PERTUZ",  "2110", 36341),
    @(40, "CC", "1143377371", "SEBASTIAN ENRIQUE VILLARREAL PERTUZ",  "2109", 36341)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
    $ws.Range("E$r").Value = $row[4]
    $ws.Range("F$r").Value = $row[5]
    $ws.Range("G$r").Value = 908526
}
